$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.006.66"
$ws.Range("E2").Value = "'  +1.08%  "
$ws.Range("D3").Value = "'3.141.73"
$ws.Range("E3").Value = "'  +1.74%  "
$ws.Range("D5").Value = "'591.21"
$ws.Range("E5").Value = "'  +1.64%  "
$ws.Range("D6").Value = "'146.25"
$ws.Range("E6").Value = "'  +1.37%  "
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("D8").Value = "'3.134.22"
$ws.Range("E8").Value = "'  +1.77%  "
$ws.Range("E9").Value = "'  +0.59%  "
$ws.Range("E10").Value = "'  +2.94%  "
$ws.Range("D11").Value = "'5.93"
$ws.Range("E11").Value = "'  +5.54%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "'  +0.43%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "'  +1.33%  "
$ws.Range("D14").Value = "'37.20"
$ws.Range("E14").Value = "'  -1.72%  "
$ws.Range("D15").Value = "'3.661.34"
$ws.Range("E15").Value = "'  +1.71%  "
$ws.Range("E16").Value = "'  -0.17%  "
$ws.Range("D17").Value = "'7.26"
$ws.Range("E17").Value = "'  +2.36%  "
$ws.Range("D18").Value = "'63.801.66"
$ws.Range("E18").Value = "'  +0.92%  "
$ws.Range("D19").Value = "'3.135.60"
$ws.Range("E19").Value = "'  +1.60%  "
$ws.Range("D20").Value = "'466.37"
$ws.Range("E20").Value = "'  +1.43%  "
$ws.Range("D21").Value = "'14.35"
$ws.Range("E21").Value = "'  +1.54%  "
$ws.Range("D22").Value = "'0.731"
$ws.Range("E22").Value = "'  +1.20%  "
$ws.Range("D23").Value = "'7.58"
$ws.Range("E23").Value = "'  +1.93%  "
$ws.Range("E24").Value = "'  +12.95%  "
$ws.Range("E25").Value = "'  +1.62%  "
$ws.Range("D26").Value = "'80.91"
$ws.Range("E26").Value = "'  -0.14%  "
$ws.Range("E27").Value = "'  +0.05%  "
$ws.Range("D28").Value = "'9.82"
$ws.Range("E28").Value = "'  +10.19%  "
$ws.Range("E29").Value = "'  +1.86%  "
$ws.Range("D30").Value = "'7.35"
$ws.Range("E30").Value = "'  +8.00%  "
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("E32").Value = "'  +0.78%  "
$ws.Range("E33").Value = "'  +5.12%  "
$ws.Range("D34").Value = "'27.63"
$ws.Range("E34").Value = "'  +3.98%  "
$ws.Range("D35").Value = "'0.0₃0863"
$ws.Range("E35").Value = "'  +2.39%  "
$ws.Range("E36").Value = "'  +3.32%  "
$ws.Range("D37").Value = "'6.15"
$ws.Range("E37").Value = "'  +3.04%  "
$ws.Range("E38").Value = "'  -1.26%  "
$ws.Range("D39").Value = "'3.28"
$ws.Range("E39").Value = "'  -1.57%  "
$ws.Range("D40").Value = "'462.25"
$ws.Range("E40").Value = "'  +6.34%  "
$ws.Range("D42").Value = "'51.34"
$ws.Range("E42").Value = "'  +2.26%  "
$ws.Range("D43").Value = "'0.293"
$ws.Range("E43").Value = "'  +9.50%  "
$ws.Range("D44").Value = "'0.0372"
$ws.Range("E44").Value = "'  +1.11%  "
$ws.Range("D45").Value = "'2.890.37"
$ws.Range("E45").Value = "'  +1.19%  "
$ws.Range("D46").Value = "'40.21"
$ws.Range("E46").Value = "'  +11.79%  "
$ws.Range("E47").Value = "'  -0.50%  "
$ws.Range("D48").Value = "'132.87"
$ws.Range("E48").Value = "'  +7.23%  "
$ws.Range("D50").Value = "'0.111"
$ws.Range("E50").Value = "'  +0.90%  "
$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "'  +4.27%  "
